$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Material Type"
$ws.Range("G2:G24").Value = "DNA:Genomic"

$headerRange = $ws.Range("G1")
$headerRange.Interior.Color = 0
$headerRange.Font.Color = 16777215
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "MS Sans Serif"
$headerRange.Font.Size = 10
$headerRange.HorizontalAlignment = -4108

$dataRange = $ws.Range("G2:G24")
$dataRange.HorizontalAlignment = -4108

$ws.Range("G1").Select()
